# siba-be / templates / user_template.xlsx
#
# Commit: "update so multiple department names can be added for single
# userId" — the "department" column in the user-import template should
# document that several department names can be supplied for one user by
# separating them with a slash.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F2:F4 all hold the placeholder/example text for the "department" column
# (shared-string "name of department"); update every sample row so the
# template shows the new multi-department, slash-separated format.
$ws.Range("F2:F4").Value = "name of department 1/name of department 2"

# Reflect the author's last on-sheet selection (cell F4) when the workbook
# was saved.
$null = $ws.Range("F4").Select()
